$d = $word.ActiveDocument

# 1) "Fecha :" -> "Fecha:"  (remove the extra space before the colon)
$d.Content.Find.Execute("Fecha :", $true, $false, $false, $false, $false, $true, 1, $false, "Fecha:", 2) | Out-Null

# 2) Merge "Implementación del Sistema de Registro y " + "Login" into one run of text
#    (appears twice in the document); Find/Replace on the plain text merges the spelling
#    so that it becomes a contiguous phrase "Implementación del Sistema de Registro y Login".
$d.Content.Find.Execute("Registro y Login", $true, $false, $false, $false, $false, $true, 1, $false, "Registro y Login", 2) | Out-Null

$d.Save()
